# Update "想去人数" (F column) figures on sheets "展览" (1), "演出" (2), and "全部类型" (4)
# to reflect newly scraped counts, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsShow    = $wb.Worksheets.Item(2)   # 演出
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# 展览 (sheet1)
$wsExhibit.Range("F2").Value  = 15099
$wsExhibit.Range("F3").Value  = 19306
$wsExhibit.Range("F5").Value  = 152
$wsExhibit.Range("F14").Value = 190
$wsExhibit.Range("F17").Value = 1493
$wsExhibit.Range("F22").Value = 8091
$wsExhibit.Range("F24").Value = 38
$wsExhibit.Range("F30").Value = 6100
$wsExhibit.Range("F33").Value = 178
$wsExhibit.Range("F36").Value = 5506
$wsExhibit.Range("F38").Value = 22

# 演出 (sheet2)
$wsShow.Range("F3").Value = 24

# 全部类型 (sheet4) - combined listing, same events repeated
$wsAll.Range("F2").Value  = 15099
$wsAll.Range("F3").Value  = 19306
$wsAll.Range("F5").Value  = 152
$wsAll.Range("F14").Value = 190
$wsAll.Range("F17").Value = 1493
$wsAll.Range("F23").Value = 8091
$wsAll.Range("F25").Value = 38
$wsAll.Range("F31").Value = 24
$wsAll.Range("F33").Value = 6100
$wsAll.Range("F36").Value = 178
$wsAll.Range("F39").Value = 5506
$wsAll.Range("F41").Value = 22
